$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(53, 1).Value = 51
$ws.Cells.Item(53, 2).Value = 26577691
$ws.Cells.Item(53, 3).Value = 15838
$ws.Cells.Item(53, 4).Value = 14999.67
$ws.Cells.Item(53, 5).Value = 44278
